$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-duration"

# Version: 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
# The "Fixed Value" cell for Extension.url (row 5, column Q) repeats the
# same URL string, so it needs to track the same change.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-duration"
